# 20200427 part 2 update
# - Sheet "测试概况": fill in the CRM order numbers (column D) for the 5
#   existing test-scenario rows, and tighten the row heights back down now
#   that the long "备注" column is no longer wrapped across extra lines.
# - Sheet "BUG汇总": add the first BUG record (row 2) with all of its
#   fields, and best-fit the two date columns (D, I) to fit the new data.

$wb = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item(1)   # 测试概况
$wsBug      = $wb.Worksheets.Item(2)   # BUG汇总

# ---------------------------------------------------------------------
# Sheet 1: 测试概况 — column D (CRM order number) + row heights
# ---------------------------------------------------------------------
$wsOverview.Cells.Item(2, 4).Value = "WMX2020042705390546"
$wsOverview.Cells.Item(3, 4).Value = "WMX2020042705390547"
$wsOverview.Cells.Item(4, 4).Value = "WMX2020042705390821"
$wsOverview.Cells.Item(5, 4).Value = "WMX2020042705390832"
$wsOverview.Cells.Item(6, 4).Value = "WMX2020042705390877"

$wsOverview.Rows.Item(2).RowHeight = 24
$wsOverview.Rows.Item(3).RowHeight = 39.95
$wsOverview.Rows.Item(4).RowHeight = 24
$wsOverview.Rows.Item(5).RowHeight = 36
$wsOverview.Rows.Item(6).RowHeight = 36

# ---------------------------------------------------------------------
# Sheet 2: BUG汇总 — row 2 data
# ---------------------------------------------------------------------
$wsBug.Cells.Item(2, 1).Value = 1
$wsBug.Cells.Item(2, 2).Value = "预付费鹏博士宽带新装"
$wsBug.Cells.Item(2, 3).Value = "WMX2020042605389540"
$wsBug.Cells.Item(2, 4).Value = 43947.465277777781
$wsBug.Cells.Item(2, 5).Value = "订单从CRM受理提交到IBP后，从SOP发送订单到P7系统未成功"
$wsBug.Cells.Item(2, 6).Value = "IBP"
$wsBug.Cells.Item(2, 7).Value = "已解决"
$wsBug.Cells.Item(2, 7).Font.Color = 255
$wsBug.Cells.Item(2, 9).Value = 43947.579861111109
$wsBug.Cells.Item(2, 10).Value = "环境问题"

$wsBug.Rows.Item(2).RowHeight = 36

# Best-fit the BUG发现/解决时间 columns (D, I) now that they hold data.
$wsBug.Columns.Item(4).ColumnWidth = 16.875
$wsBug.Columns.Item(9).ColumnWidth = 16.875
